$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-2.838432111238169, -4.14933462754678, 1.666932080921356),
    @(-4.300542519280826, 2.241362005864094, 4.135142751430305),
    @(-1.90287771350459, 13.77420247228532, -1.556544103120493),
    @(-0.1651651027956011, 2.663977503776609, -2.861205784898042),
    @(8.837720880383007, -7.873765813676991, -5.242947754107039),
    @(0.8822141462250848, -0.9568282867732949, -2.844946472268377),
    @(-2.166194288354044, 0.3768198490145664, 1.521733685543742),
    @(-9.005332197798026, -16.56740246321032, -2.677433475067815),
    @(0.6589992062041609, -9.836473477514197, 4.410142396625722),
    @(1.44424632505363, -3.936239650375036, 3.583756440172107)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
